$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.429.65"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.942.87"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.361"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0843"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "2.228.02"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.811"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("D17").Value = "1.945.40"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "36.384.86"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("E24").Value = "  -5.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("E28").Value = "  +4.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.51%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -5.02%  "
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("B34").Value = "THORChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.90%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.21%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0988"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").Value = "1.341.18"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "2.119.25"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.25%  "
